# Weekly refresh of the "Camote" price series: the daily logic shifted one
# row back in the historical window, dropping the oldest record and
# appending a new week (row 10) at the end. Column D (Fecha) and the price
# columns J/K/L/M/P are the ones that move; the descriptive columns
# (A,B,C,E,F,G,H,I,N,O,Q,R) are constant for every row of this subset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for Fecha (D), Volumen (J), Precio minimo (K), Precio maximo (L),
# Precio promedio ponderado (M) and Precio $/Kg (P) for rows 2..10.
$rows = @{
    2  = @(44792, 160, 9000,  10000, 9500,  528)
    3  = @(44804, 50,  9500,  10000, 9750,  542)
    4  = @(45205, 200, 11000, 12000, 11500, 639)
    5  = @(45092, 210, 10000, 11000, 10714, 595)
    6  = @(45273, 200, 9000,  10000, 9500,  528)
    7  = @(45245, 100, 9000,  10000, 9500,  528)
    8  = @(44714, 80,  9000,  10000, 9500,  528)
    9  = @(45175, 250, 11000, 12000, 11500, 639)
    10 = @(45215, 200, 11000, 12000, 11500, 639)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]

    $ws.Cells.Item($r, 4).Value  = $vals[0]   # D - Fecha
    $ws.Cells.Item($r, 10).Value = $vals[1]   # J - Volumen
    $ws.Cells.Item($r, 11).Value = $vals[2]   # K - Precio minimo
    $ws.Cells.Item($r, 12).Value = $vals[3]   # L - Precio maximo
    $ws.Cells.Item($r, 13).Value = $vals[4]   # M - Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $vals[5]   # P - Precio $/Kg
}

# Row 10 is brand new - fill in the columns that stay constant across the
# whole sheet, copied from row 9, and match its date-cell number format.
# NOTE: reading `.Value` needs explicit call syntax `()` in this host, else
# it yields the property descriptor instead of invoking the getter.
$ws.Cells.Item(10, 1).Value  = $ws.Cells.Item(9, 1).Value()    # A - Mercado ID
$ws.Cells.Item(10, 2).Value  = $ws.Cells.Item(9, 2).Value()    # B - Mercado
$ws.Cells.Item(10, 3).Value  = $ws.Cells.Item(9, 3).Value()    # C - Region
$ws.Cells.Item(10, 5).Value  = $ws.Cells.Item(9, 5).Value()    # E - Codreg
$ws.Cells.Item(10, 6).Value  = $ws.Cells.Item(9, 6).Value()    # F - Categoria ID
$ws.Cells.Item(10, 7).Value  = $ws.Cells.Item(9, 7).Value()    # G - Categoria
$ws.Cells.Item(10, 8).Value  = $ws.Cells.Item(9, 8).Value()    # H - Variedad
$ws.Cells.Item(10, 9).Value  = $ws.Cells.Item(9, 9).Value()    # I - Calidad
$ws.Cells.Item(10, 14).Value = $ws.Cells.Item(9, 14).Value()   # N - Unidad de comercializacion
$ws.Cells.Item(10, 15).Value = $ws.Cells.Item(9, 15).Value()   # O - Origen
$ws.Cells.Item(10, 17).Value = $ws.Cells.Item(9, 17).Value()   # Q - Kg o Unidades
$ws.Cells.Item(10, 18).Value = $ws.Cells.Item(9, 18).Value()   # R - Clasificacion

$ws.Cells.Item(10, 4).NumberFormat = $ws.Cells.Item(9, 4).NumberFormat  # D - Fecha format
